$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells for Apple/Bowl feature-value columns to lowercase
$ws.Range("B2").Value2 = "apple_feature"
$ws.Range("C2").Value2 = "apple_val"
$ws.Range("D2").Value2 = "bowl_feature"
$ws.Range("E2").Value2 = "bowl_val"

# Rename the "Concepts" column entries for Apple/Bowl to lowercase
$ws.Range("A3").Value2 = "apple"
$ws.Range("A4").Value2 = "bowl"

# Add the missing "spoon" concept label in row 11
$ws.Range("A11").Value2 = "spoon"

# Update the active selection: drop the C1 scroll anchor and select E3
[void]$ws.Range("E3").Select()
